$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 2).Value = 0
$ws.Cells.Item(1, 3).Value = 1
$ws.Cells.Item(1, 4).Value = 2
$ws.Cells.Item(1, 5).Value = 3
$ws.Cells.Item(1, 6).Value = 4
$ws.Cells.Item(1, 7).Value = 5
$ws.Cells.Item(1, 8).Value = 6
$ws.Cells.Item(1, 9).Value = 7
$ws.Cells.Item(1, 10).Value = 8
$ws.Cells.Item(1, 11).Value = 9
$ws.Cells.Item(1, 12).Value = 10
$ws.Cells.Item(1, 13).Value = 11
$ws.Cells.Item(1, 14).Value = 12
$ws.Cells.Item(1, 15).Value = 13
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15
$ws.Cells.Item(1, 18).Value = 16
$ws.Cells.Item(1, 19).Value = 17
$ws.Cells.Item(1, 20).Value = 18
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "HKL"
$ws.Cells.Item(2, 3).Value = "[1, 1, 0]"
$ws.Cells.Item(2, 4).Value = "[2, 0, 0]"
$ws.Cells.Item(2, 5).Value = "[2, 1, 1]"
$ws.Cells.Item(2, 6).Value = "[2, 2, 0]"
$ws.Cells.Item(2, 7).Value = "[3, 1, 0]"
$ws.Cells.Item(2, 8).Value = "[2, 2, 2]"
$ws.Cells.Item(2, 9).Value = "[3, 2, 1]"
$ws.Cells.Item(2, 10).Value = "[4, 0, 0]"
$ws.Cells.Item(2, 11).Value = "1Pair-A"
$ws.Cells.Item(2, 12).Value = "1Pair-B"
$ws.Cells.Item(2, 13).Value = "2Pairs-A"
$ws.Cells.Item(2, 14).Value = "2Pairs-B"
$ws.Cells.Item(2, 15).Value = "3Pairs-A"
$ws.Cells.Item(2, 16).Value = "3Pairs-B"
$ws.Cells.Item(2, 17).Value = "3Pairs-C"
$ws.Cells.Item(2, 18).Value = "4Pairs"
$ws.Cells.Item(2, 19).Value = "5A4F"
$ws.Cells.Item(2, 20).Value = "MaxUnique"
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "Equal Angle"
$ws.Cells.Item(3, 3).Value = 0.9731772334293948
$ws.Cells.Item(3, 4).Value = 0.8101512968299712
$ws.Cells.Item(3, 5).Value = 1.077168587896254
$ws.Cells.Item(3, 6).Value = 0.9731772334293948
$ws.Cells.Item(3, 7).Value = 0.8578386167146974
$ws.Cells.Item(3, 8).Value = 1.266469740634006
$ws.Cells.Item(3, 9).Value = 1.052348703170029
$ws.Cells.Item(3, 10).Value = 0.8101512968299712
$ws.Cells.Item(3, 11).Value = 0.9731772334293948
$ws.Cells.Item(3, 12).Value = 1.077168587896254
$ws.Cells.Item(3, 13).Value = 0.9436599423631125
$ws.Cells.Item(3, 14).Value = 0.9436599423631125
$ws.Cells.Item(3, 15).Value = 0.9150528338136409
$ws.Cells.Item(3, 16).Value = 0.9534990393852065
$ws.Cells.Item(3, 17).Value = 0.9534990393852065
$ws.Cells.Item(3, 18).Value = 0.9584185878962536
$ws.Cells.Item(3, 19).Value = 0.9584185878962536
$ws.Cells.Item(3, 20).Value = 1.006192363112392
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "CLR"
$ws.Cells.Item(4, 3).Value = 1.008363564545268
$ws.Cells.Item(4, 4).Value = 0.9671685740387875
$ws.Cells.Item(4, 5).Value = 0.9986341983137771
$ws.Cells.Item(4, 6).Value = 1.008363564545268
$ws.Cells.Item(4, 7).Value = 0.9773359313944642
$ws.Cells.Item(4, 8).Value = 1.003080432085543
$ws.Cells.Item(4, 9).Value = 1.001558435594122
$ws.Cells.Item(4, 10).Value = 0.9671685740387875
$ws.Cells.Item(4, 11).Value = 1.008363564545268
$ws.Cells.Item(4, 12).Value = 0.9986341983137771
$ws.Cells.Item(4, 13).Value = 0.9829013861762823
$ws.Cells.Item(4, 14).Value = 0.9829013861762823
$ws.Cells.Item(4, 15).Value = 0.981046234582343
$ws.Cells.Item(4, 16).Value = 0.9913887789659442
$ws.Cells.Item(4, 17).Value = 0.9913887789659442
$ws.Cells.Item(4, 18).Value = 0.9956324753607753
$ws.Cells.Item(4, 19).Value = 0.9956324753607753
$ws.Cells.Item(4, 20).Value = 0.9926901893286603
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "BT8Hex"
$ws.Cells.Item(5, 3).Value = 1.023469364379357
$ws.Cells.Item(5, 4).Value = 0.9436377789349696
$ws.Cells.Item(5, 5).Value = 1.000695745163666
$ws.Cells.Item(5, 6).Value = 1.023469364379357
$ws.Cells.Item(5, 7).Value = 0.9655219032205984
$ws.Cells.Item(5, 8).Value = 1.007429030350542
$ws.Cells.Item(5, 9).Value = 1.007177195907384
$ws.Cells.Item(5, 10).Value = 0.9436377789349696
$ws.Cells.Item(5, 11).Value = 1.023469364379357
$ws.Cells.Item(5, 12).Value = 1.000695745163666
$ws.Cells.Item(5, 13).Value = 0.9721667620493177
$ws.Cells.Item(5, 14).Value = 0.9721667620493177
$ws.Cells.Item(5, 15).Value = 0.9699518091064112
$ws.Cells.Item(5, 16).Value = 0.9892676294926641
$ws.Cells.Item(5, 17).Value = 0.9892676294926641
$ws.Cells.Item(5, 18).Value = 0.9978180632143374
$ws.Cells.Item(5, 19).Value = 0.9978180632143374
$ws.Cells.Item(5, 20).Value = 0.9913218363260862
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "Spiral"
$ws.Cells.Item(6, 3).Value = 0.9901352758028489
$ws.Cells.Item(6, 4).Value = 0.9935058341818552
$ws.Cells.Item(6, 5).Value = 0.9974273646219645
$ws.Cells.Item(6, 6).Value = 0.9901352758028489
$ws.Cells.Item(6, 7).Value = 0.9896182599373986
$ws.Cells.Item(6, 8).Value = 1.00343188788881
$ws.Cells.Item(6, 9).Value = 0.9958916829036506
$ws.Cells.Item(6, 10).Value = 0.9935058341818552
$ws.Cells.Item(6, 11).Value = 0.9901352758028489
$ws.Cells.Item(6, 12).Value = 0.9974273646219645
$ws.Cells.Item(6, 13).Value = 0.9954665994019098
$ws.Cells.Item(6, 14).Value = 0.9954665994019098
$ws.Cells.Item(6, 15).Value = 0.9935171529137395
$ws.Cells.Item(6, 16).Value = 0.9936894915355562
$ws.Cells.Item(6, 17).Value = 0.9936894915355561
$ws.Cells.Item(6, 18).Value = 0.9928009376023793
$ws.Cells.Item(6, 19).Value = 0.9928009376023793
$ws.Cells.Item(6, 20).Value = 0.995001717556088
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "OffsetF"
$ws.Cells.Item(7, 3).Value = 1.121342149089212
$ws.Cells.Item(7, 4).Value = 1.094052022585041
$ws.Cells.Item(7, 5).Value = 0.8937743216543073
$ws.Cells.Item(7, 6).Value = 1.121342149089212
$ws.Cells.Item(7, 7).Value = 1.092233990538679
$ws.Cells.Item(7, 8).Value = 0.6889417584802076
$ws.Cells.Item(7, 9).Value = 0.9523938678260909
$ws.Cells.Item(7, 10).Value = 1.094052022585041
$ws.Cells.Item(7, 11).Value = 1.121342149089212
$ws.Cells.Item(7, 12).Value = 0.8937743216543073
$ws.Cells.Item(7, 13).Value = 0.9939131721196739
$ws.Cells.Item(7, 14).Value = 0.9939131721196739
$ws.Cells.Item(7, 15).Value = 1.026686778259342
$ws.Cells.Item(7, 16).Value = 1.036389497776186
$ws.Cells.Item(7, 17).Value = 1.036389497776186
$ws.Cells.Item(7, 18).Value = 1.057627660604443
$ws.Cells.Item(7, 19).Value = 1.057627660604443
$ws.Cells.Item(7, 20).Value = 0.9737896850289228
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "OffsetA"
$ws.Cells.Item(8, 3).Value = 0.902126487668975
$ws.Cells.Item(8, 4).Value = 1.070736575035974
$ws.Cells.Item(8, 5).Value = 1.014757552172515
$ws.Cells.Item(8, 6).Value = 0.902126487668975
$ws.Cells.Item(8, 7).Value = 1.024146020964109
$ws.Cells.Item(8, 8).Value = 1.035252338869765
$ws.Cells.Item(8, 9).Value = 0.9774295127375004
$ws.Cells.Item(8, 10).Value = 1.070736575035974
$ws.Cells.Item(8, 11).Value = 0.902126487668975
$ws.Cells.Item(8, 12).Value = 1.014757552172515
$ws.Cells.Item(8, 13).Value = 1.042747063604244
$ws.Cells.Item(8, 14).Value = 1.042747063604244
$ws.Cells.Item(8, 15).Value = 1.036546716057533
$ws.Cells.Item(8, 16).Value = 0.9958735382924879
$ws.Cells.Item(8, 17).Value = 0.995873538292488
$ws.Cells.Item(8, 18).Value = 0.9724367756366098
$ws.Cells.Item(8, 19).Value = 0.9724367756366098
$ws.Cells.Item(8, 20).Value = 1.00407474790814
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "RD Single"
$ws.Cells.Item(9, 3).Value = 0.84
$ws.Cells.Item(9, 4).Value = 0.44
$ws.Cells.Item(9, 5).Value = 1.27
$ws.Cells.Item(9, 6).Value = 0.84
$ws.Cells.Item(9, 7).Value = 0.68
$ws.Cells.Item(9, 8).Value = 1.52
$ws.Cells.Item(9, 9).Value = 1.15
$ws.Cells.Item(9, 10).Value = 0.44
$ws.Cells.Item(9, 11).Value = 0.84
$ws.Cells.Item(9, 12).Value = 1.27
$ws.Cells.Item(9, 13).Value = 0.855
$ws.Cells.Item(9, 14).Value = 0.855
$ws.Cells.Item(9, 15).Value = 0.7966666666666667
$ws.Cells.Item(9, 16).Value = 0.85
$ws.Cells.Item(9, 17).Value = 0.85
$ws.Cells.Item(9, 18).Value = 0.8474999999999999
$ws.Cells.Item(9, 19).Value = 0.8474999999999999
$ws.Cells.Item(9, 20).Value = 0.9833333333333334
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "TD Single"
$ws.Cells.Item(10, 3).Value = 1.98
$ws.Cells.Item(10, 4).Value = 0.21
$ws.Cells.Item(10, 5).Value = 0.84
$ws.Cells.Item(10, 6).Value = 1.98
$ws.Cells.Item(10, 7).Value = 0.64
$ws.Cells.Item(10, 8).Value = 0.69
$ws.Cells.Item(10, 9).Value = 1.14
$ws.Cells.Item(10, 10).Value = 0.21
$ws.Cells.Item(10, 11).Value = 1.98
$ws.Cells.Item(10, 12).Value = 0.84
$ws.Cells.Item(10, 13).Value = 0.525
$ws.Cells.Item(10, 14).Value = 0.525
$ws.Cells.Item(10, 15).Value = 0.5633333333333334
$ws.Cells.Item(10, 16).Value = 1.01
$ws.Cells.Item(10, 17).Value = 1.01
$ws.Cells.Item(10, 18).Value = 1.2525
$ws.Cells.Item(10, 19).Value = 1.2525
$ws.Cells.Item(10, 20).Value = 0.9166666666666665
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(11, 3).Value = 0.9884187614602742
$ws.Cells.Item(11, 4).Value = 0.9966049724701519
$ws.Cells.Item(11, 5).Value = 0.9969555234280741
$ws.Cells.Item(11, 6).Value = 0.9884187614602742
$ws.Cells.Item(11, 7).Value = 0.9919241364401797
$ws.Cells.Item(11, 8).Value = 0.9998624821959811
$ws.Cells.Item(11, 9).Value = 0.994928375453311
$ws.Cells.Item(11, 10).Value = 0.9966049724701519
$ws.Cells.Item(11, 11).Value = 0.9884187614602742
$ws.Cells.Item(11, 12).Value = 0.9969555234280741
$ws.Cells.Item(11, 13).Value = 0.996780247949113
$ws.Cells.Item(11, 14).Value = 0.996780247949113
$ws.Cells.Item(11, 15).Value = 0.9951615441128019
$ws.Cells.Item(11, 16).Value = 0.9939930857861667
$ws.Cells.Item(11, 17).Value = 0.9939930857861667
$ws.Cells.Item(11, 18).Value = 0.9925995047046936
$ws.Cells.Item(11, 19).Value = 0.9925995047046936
$ws.Cells.Item(11, 20).Value = 0.9947823752413286

# Apply the bold/centered/bordered header style (matching the existing header row/column)
# to the newly introduced header cells by copying formatting from an existing styled neighbor.
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)

$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)
